# Week 7 backlog clean-up + new "Week8" section
# ------------------------------------------------------------------
# - "Finish the pun engine" gets split in two (so the existing
#   _GoBack bookmark sits mid-run) and highlighted yellow.
# - "….Find a new objective off of the proposal" gets highlighted
#   green and loses the bookmark (it now lives in the paragraph
#   above).
# - The old empty "Optional Objectives"-style (italic+underline)
#   paragraph becomes the new "Week8" heading (underline only).
# - Two new backlog bullets are appended under Week8.

$d = $word.ActiveDocument

function Get-ParaIndexByText($doc, $text) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        $rng = $doc.Range($p.Range.Start, $p.Range.End)
        $t = $rng.Text.TrimEnd([char]13, [char]7)
        if ($t -eq $text) {
            return $i
        }
    }
    return -1
}

# Locate the "Finish the pun engine" paragraph; the two paragraphs
# that follow it are the objective line and the (empty) italic/
# underline spacer paragraph that precedes "Week 7"'s successor.
$startIndex = Get-ParaIndexByText $d "Finish the pun engine"
if ($startIndex -eq -1) {
    throw "Could not find anchor paragraph 'Finish the pun engine'"
}

$firstPara = $d.Paragraphs.Item($startIndex)
$lastPara = $d.Paragraphs.Item($startIndex + 2)

$target = $d.Range($firstPara.Range.Start, $lastPara.Range.End)

$newXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t>Finish th</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t>e pun engine</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:highlight w:val="green"/></w:rPr><w:t>&#8230;.Find a new objective off of the proposal</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>Week8</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Tweak the Pun engine (more faster, more ways to find a noun/</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>verb</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>,etc</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:t>)</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Add more domains to the Library for the Pun Engine.</w:t></w:r></w:p>
'@

[void]$target.InsertXML($newXml)

Write-Output "Backlog Week 7/Week8 section updated."
